# Remove unused columns from the list report test file.
# Original layout: A=Created On, B=Created By, C=Changed On, D=Changed By,
#                   E=OrderNo, F=Customer (hyperlinked "mailto" email), G=Currency
# Target layout:    A=OrderNo, B=Customer (hyperlinked)
# i.e. delete columns A:D (shifting E->A, F->B, G->C), then delete the now
# orphaned Currency column (new C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the hyperlink target so we can re-create the links after the
# column shuffle (the Hyperlinks collection does not auto-track cells that
# move via a column delete).
$linkAddress = "mailto:test@test.de"

$ws.Range("A1:D1").EntireColumn.Delete()
$ws.Range("C1").EntireColumn.Delete()

# Re-point the two hyperlinks at their new home (column B) and restore the
# "Link" cell style that the column deletes already carried along with the
# cell contents.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $linkAddress)
$ws.Hyperlinks.Add($ws.Range("B3"), $linkAddress)
$ws.Range("B2").Style = "Link"
$ws.Range("B3").Style = "Link"

# Reflect the user's final selection (columns A:D were selected/highlighted
# right before saving).
$ws.Range("A1:D1048576").Select()
